# Swap the "Integral" theme colours on theme1.xml (the deck's real, reachable
# theme part, linked from the slide master) for the "Office Theme" colour
# scheme. The font scheme and format scheme are already byte-identical
# between the old theme1.xml ("Integral") and theme2.xml ("Office Theme"),
# so only the 12 colour-scheme slots need to change.
#
# ThemeColorScheme.Item(n) order is: dk1, lt1, dk2, lt2, accent1..accent6,
# hlink, folHlink. The .RGB property takes/returns a packed BGR integer
# (same convention as VBA's RGB()/OLE_COLOR), i.e. 0xBBGGRR.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$tcs = $m.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388     # dk2      44546A
$tcs.Item(4).RGB  = 15132391    # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939    # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501     # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845    # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407       # accent4  FFC000
$tcs.Item(9).RGB  = 12874308    # accent5  4472C4
$tcs.Item(10).RGB = 4697456     # accent6  70AD47
$tcs.Item(11).RGB = 12673797    # hlink    0563C1
$tcs.Item(12).RGB = 7491477     # folHlink 954F72
